# Fix target related bugs
# - Rows 2-21 get new BRAND / Item Name / UOM values (new GPM lineup),
#   refreshed BSL NO, BE (MTD Sales Target Value) and BK (GPMID) figures,
#   and the salesman name in BL is updated.
# - The former row 22 (Softi / BOG GPM) is removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$salesman = "Mr. Md. Tafsir Bashar"

$rows = @(
    @{ Row=2;  A=4;   B="Aldorin";  C="Aldorin 50mg Tablet - 24's";          D="24's";  BE=179.91; BK=3216 },
    @{ Row=3;  A=17;  B="Cardoneb"; C="Cardoneb 5 FC Tablet 30's";           D="30's";  BE=224.89; BK=3216 },
    @{ Row=4;  A=17;  B="Cardoneb"; C="Cardoneb 2.5 FC Tablet 30's";         D="30's";  BE=157.42; BK=3216 },
    @{ Row=5;  A=18;  B="Cardovan"; C="Cardovan Plus 80/12.5 Tablet 30's";   D="30's";  BE=224.89; BK=3216 },
    @{ Row=6;  A=18;  B="Cardovan"; C="Cardovan 160mg Tablet 30's";         D="30's";  BE=359.82; BK=3216 },
    @{ Row=7;  A=18;  B="Cardovan"; C="Cardovan 80mg Tablet 30's";          D="30's";  BE=224.89; BK=3216 },
    @{ Row=8;  A=18;  B="Cardovan"; C="Cardovan Plus 160/12.5 Tablet 30's"; D="30's";  BE=359.82; BK=3216 },
    @{ Row=9;  A=26;  B="Dialon";   C="Dialon 4mg Tablet";                  D="20'S";  BE=179.91; BK=3216 },
    @{ Row=10; A=52;  B="GLIKAZID"; C="Glikazid 80mg Tablet 30's";          D="30's";  BE=112.44; BK=3216 },
    @{ Row=11; A=57;  B="Irbes";    C="Irbes 75mg Tablet";                  D="50 's"; BE=187.41; BK=3216 },
    @{ Row=12; A=68;  B="Ligazid";  C="Ligazid 5mg Tablet 20's";            D="20's";  BE=224.89; BK=3216 },
    @{ Row=13; A=68;  B="Ligazid";  C="Ligazid 5mg Tablet 10's";            D="10's";  BE=112.44; BK=3216 },
    @{ Row=14; A=68;  B="Ligazid";  C="Ligazid M 2.5/500";                  D="20's";  BE=179.91; BK=3216 },
    @{ Row=15; A=70;  B="Lipicon";  C="Lipicon 20mg Tablet - 20's";         D="20 's"; BE=269.87; BK=3216 },
    @{ Row=16; A=70;  B="Lipicon";  C="Lipicon 10mg Tablet Container 30's"; D="30's";  BE=224.89; BK=3216 },
    @{ Row=17; A=70;  B="Lipicon";  C="Lipicon 40mg Tablet - 10's";         D="10 's"; BE=179.91; BK=3216 },
    @{ Row=18; A=70;  B="Lipicon";  C="Lipicon 10mg Tablet - 40's";         D="40 's"; BE=299.85; BK=3216 },
    @{ Row=19; A=104; B="Pivasta";  C="Pivasta 2mg Tablet 20's";            D="20's";  BE=149.93; BK=3216 },
    @{ Row=20; A=123; B="Sitazid";  C="Sitazid 50mg Tablet 20's";           D="20's";  BE=194.90; BK=3216 },
    @{ Row=21; A=123; B="Sitazid";  C="Sitazid 100mg Tablet 10's";          D="10's";  BE=187.41; BK=3216 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.A
    $ws.Range("B$n").Value = $r.B
    $ws.Range("C$n").Value = $r.C
    $ws.Range("D$n").Value = $r.D
    $ws.Range("BE$n").Value = $r.BE
    $ws.Range("BK$n").Value = $r.BK
    $ws.Range("BL$n").Value = $salesman
}

# The old row 22 (BSL NO 165 / Softi / BOG) no longer exists in the refreshed target list.
$ws.Rows.Item(22).Delete()
